# setDataIntoExcel method created in Excel Utility
# This script reproduces a data-write operation performed against the
# "prod" worksheet of the workbook (as exercised by the Selenium/Java
# Excel utility's setDataIntoExcel method), and leaves that sheet as the
# active tab, the way Excel records it after a save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prod")

# Overwrite the previous "oneplus" row with a marker value used by the
# new utility method, and drop the now-unused second column value.
$ws.Range("A4").Value = "SetValue-Working"
$ws.Range("B4").ClearContents()

# New data points written further down the sheet.
$ws.Range("A6").Value = "SetValue-Working"
$ws.Range("B7").Value = "SetValue-Working"

# Make "prod" the active/selected sheet (previously "contact" was active).
$ws.Activate()
